$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Menu List" sheet becomes the active sheet/tab (workbook activeTab 3 -> 2),
#    its selection moves to G33, and it gains tabSelected="1" while the
#    previously active "Modifier List" sheet loses it automatically.
# ---------------------------------------------------------------------------
$menuList = $wb.Worksheets.Item("Menu List")
$menuList.Activate()
$menuList.Range("G33").Select()

# ---------------------------------------------------------------------------
# 2. Menu List rows 22-24: refresh the random identifier text in column A
#    (same value reused on all three rows) and update the "Required/Min/Max"
#    quantity markers on row 23 from 1/4 to 3/3.
# ---------------------------------------------------------------------------
$menuList.Range("A22").Value = "pFiMb9HUdJ"
$menuList.Range("A23").Value = "pFiMb9HUdJ"
$menuList.Range("A24").Value = "pFiMb9HUdJ"

foreach ($addr in @("F23", "G23", "H23", "I23")) {
    $cell = $menuList.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "3"
}

# ---------------------------------------------------------------------------
# 3. Strip the now-unwanted cell-level formatting on rows 22-24:
#    - column A keeps center/top alignment with no border (column default)
#    - columns B:D keep center alignment with no border (column default)
#    - the highlighted quantity cells (F/G/H/I/J/K) revert to plain/default
# ---------------------------------------------------------------------------
$colA = $menuList.Range("A22:A24")
$colA.HorizontalAlignment = -4108   # xlCenter
$colA.VerticalAlignment = -4160     # xlTop
$colA.Borders.LineStyle = -4142     # xlLineStyleNone

$colBD = $menuList.Range("B22:D24")
$colBD.HorizontalAlignment = -4108  # xlCenter
$colBD.Borders.LineStyle = -4142    # xlLineStyleNone

$menuList.Range("F22:G22").Style = "Normal"
$menuList.Range("F23:I23").Style = "Normal"
$menuList.Range("F24:K24").Style = "Normal"

# ---------------------------------------------------------------------------
# 4. "Modifier List" sheet selection no longer carries the active tab marker;
#    nothing else on this sheet changes besides losing tabSelected (handled
#    automatically above by activating "Menu List" instead).
# ---------------------------------------------------------------------------
